$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.961.46"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "2.272.71"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("D5").Value = "231.07"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").Value = "0.637"
$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("D7").Value = "63.62"
$ws.Range("E7").Value = "  +3.12%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +6.86%  "

$ws.Range("D10").Value = "0.0995"
$ws.Range("E10").Value = "  +7.35%  "

$ws.Range("D11").Value = "57.58"
$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("D12").Value = "27.41"
$ws.Range("E12").Value = "  +15.25%  "

$ws.Range("D14").Value = "2.613.62"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").Value = "15.80"
$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").Value = "6.13"
$ws.Range("E16").Value = "  +6.40%  "

$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  +3.30%  "

$ws.Range("D18").Value = "2.266.64"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").Value = "43.911.92"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("E20").Value = "  +7.25%  "

$ws.Range("D21").Value = "73.82"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  -2.40%  "

$ws.Range("D23").Value = "252.55"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  -4.44%  "

$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -2.12%  "

$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +2.14%  "

$ws.Range("E28").Value = "  +24.85%  "

$ws.Range("D29").Value = "171.80"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("E30").Value = "  -2.58%  "

$ws.Range("D31").Value = "20.93"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -5.80%  "

$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("E34").Value = "  +6.31%  "

$ws.Range("E35").Value = "  +0.64%  "

$ws.Range("E36").Value = "  -3.71%  "

$ws.Range("D37").Value = "3.81"
$ws.Range("E37").Value = "  +4.15%  "

$ws.Range("D38").Value = "6.54"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("E39").Value = "  -4.91%  "

$ws.Range("D40").Value = "0.0260"
$ws.Range("E40").Value = "  +3.10%  "

$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").Value = "0.000227"
$ws.Range("E42").Value = "  +3.88%  "

$ws.Range("D43").Value = "0.0990"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").Value = "17.64"
$ws.Range("E44").Value = "  +4.81%  "

$ws.Range("D45").Value = "8.28"
$ws.Range("E45").Value = "  -6.46%  "

$ws.Range("D46").Value = "10.47"
$ws.Range("E46").Value = "  +10.38%  "

$ws.Range("D47").Value = "98.36"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("D49").Value = "4.37"
$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("D50").Value = "1.446.39"
$ws.Range("E50").Value = "  -1.98%  "

$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  +0.64%  "
